$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 181
$ws1.Range("F3").Value = 103
$ws1.Range("F5").Value = 954
$ws1.Range("F6").Value = 5134
$ws1.Range("F7").Value = 419
$ws1.Range("F8").Value = 597
$ws1.Range("F11").Value = 69
$ws1.Range("F12").Value = 22
$ws1.Range("F13").Value = 550
$ws1.Range("F14").Value = 5
$ws1.Range("F15").Value = 15
$ws1.Range("F17").Value = 1694
$ws1.Range("F18").Value = 1437
$ws1.Range("F19").Value = 781
$ws1.Range("F22").Value = 287
$ws1.Range("F28").Value = 2302
$ws1.Range("F29").Value = 164
$ws1.Range("F30").Value = 88
$ws1.Range("F31").Value = 66
$ws1.Range("F33").Value = 230
$ws1.Range("F35").Value = 41
$ws1.Range("F39").Value = 604
$ws1.Range("F41").Value = 37
$ws1.Range("F43").Value = 53

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 181
$ws4.Range("F4").Value = 103
$ws4.Range("F5").Value = 954
$ws4.Range("F7").Value = 5134
$ws4.Range("F8").Value = 419
$ws4.Range("F9").Value = 597
$ws4.Range("F16").Value = 69
$ws4.Range("F17").Value = 22
$ws4.Range("F18").Value = 550
$ws4.Range("F19").Value = 5
$ws4.Range("F20").Value = 15
$ws4.Range("F23").Value = 1694
$ws4.Range("F24").Value = 1437
$ws4.Range("F25").Value = 781
$ws4.Range("F28").Value = 287
$ws4.Range("F34").Value = 2302
$ws4.Range("F35").Value = 164
$ws4.Range("F36").Value = 88
$ws4.Range("F38").Value = 230
$ws4.Range("F40").Value = 41
$ws4.Range("F43").Value = 604
$ws4.Range("F45").Value = 37
$ws4.Range("F47").Value = 53
$ws4.Range("F49").Value = 3
